# "Generate Report for Handoff"
# The file e0c92acc-1156-4d4b-a284-c032b531ec62.md has finished translation
# and is now ready for handoff. Update its status + timestamps on every
# sheet of the report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for e0c92acc-1156-4d4b-a284-c032b531ec62.md is row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-14-20 02:14:50"

# --- zh-cn sheet: same file is row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-20 02:14:47"

# --- de-de sheet: same file is row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-20 02:14:50"
